# "Final Resume4 and created Assignment5"
# Adds a "Задание" (Assignment) column, renames/repurposes the old "Примечание"
# column into "Степень выполнения" (Выполнено / Не выполнено), updates several
# students' scores, recolors the name cells to reflect degree of completion,
# and appends a new student row (Сорокин) highlighted in blue.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color helper (BGR-packed long, same as VBA's RGB())
function RGBColor([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$green      = RGBColor 0   176 80    # FF00B050 - "Выполнено" (best)
$lightGreen = RGBColor 146 208 80    # FF92D050 - "Выполнено"
$grey25     = RGBColor 191 191 191   # White, Background 1, Darker 25% (~theme0 tint -0.25)
$blue       = RGBColor 0   176 240   # FF00B0F0 - new student highlight

# ---- Header row ----
$ws.Range("C1").Value = "Степень выполнения"
$ws.Range("D1").Value = "Задание"

# ---- Row 2: Асеев ----
$ws.Range("A2").Interior.Color = $lightGreen
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = "Выполнено"
$ws.Range("D2").Value = "Назначение системы, рецензирование документации второй подгруппы"

# ---- Row 3: Акимутин ----
$ws.Range("A3").Interior.Color = $lightGreen
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "Выполнено"
$ws.Range("D3").Value = "Пользовательские требования, работа компоненты тестирования (как будет происходить проверка задач)"

# ---- Row 4: Бидзиля ----
$ws.Range("A4").Interior.Color = $green
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = "Выполнено"
$ws.Range("D4").Value = "Проектирование сущностей и связей БД, черновик календарногоплана проекта, выбор жизненного цикла"

# ---- Row 5: Бурамбекова ----
$ws.Range("A5").Interior.Color = $lightGreen
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = "Выполнено"
$ws.Range("D5").Value = "Рецензирование документации второй подгруппы"

# ---- Row 6: Заварзин ----
$ws.Range("A6").Interior.Color = $grey25
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "Не выполнено"
$ws.Range("D6").Value = "Выбор жизненного цикла"

# ---- Row 7: Лазарев ----
$ws.Range("A7").Interior.Color = $lightGreen
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = "Выполнено"
$ws.Range("D7").Value = "Написание глоссария"

# ---- Row 8: Малофеева ----
$ws.Range("A8").Interior.Color = $lightGreen
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = "Выполнено"
$ws.Range("D8").Value = "Составление новой версии требований"

# ---- Row 9: Петров ----
$ws.Range("A9").Interior.Color = $lightGreen
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = "Выполнено"
$ws.Range("D9").Value = "Цели системы, подготовка к ответу на семинаре"

# ---- Row 10: Руданов ----
$ws.Range("A10").Interior.Color = $lightGreen
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "Выполнено"
$ws.Range("D10").Value = "Выбор языка программирования и среды разработки с обоснованием, подготовка к ответу на семинаре"

# ---- Row 11: new student Сорокин ----
$ws.Range("A11").Value = "Сорокин"
$ws.Range("A11").Interior.Color = $blue
$ws.Range("B11").Value = "-"
$ws.Range("B11").HorizontalAlignment = -4152   # xlRight
$ws.Range("C11").Value = "Выполнено"
$ws.Range("D11").Value = "Черновик календарного плана проекта"

# ---- Column widths (best-fit approximation for the new/changed columns) ----
# (the engine quantizes stored width to a 1/6-character pixel grid, so these
# inputs are chosen to land on the closest reachable value to the true
# best-fit widths of 20.7109375 / 101.5703125)
$ws.Columns.Item(3).ColumnWidth = 19.75
$ws.Columns.Item(4).ColumnWidth = 100.59

# ---- Selection, matching the saved workbook state ----
[void]$ws.Range("B15").Select()
